$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and add the four new sheets in final order:
#    Comm, Property, Guild, Tip, Item
# ---------------------------------------------------------------------------
$sheetComm = $wb.Worksheets.Item(1)
$sheetComm.Name = "Comm"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetProperty = $wb.Worksheets.Add($null, $lastSheet)
$sheetProperty.Name = "Property"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetGuild = $wb.Worksheets.Add($null, $lastSheet)
$sheetGuild.Name = "Guild"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetTip = $wb.Worksheets.Add($null, $lastSheet)
$sheetTip.Name = "Tip"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetItem = $wb.Worksheets.Add($null, $lastSheet)
$sheetItem.Name = "Item"

# ---------------------------------------------------------------------------
# 2. "Comm" sheet (was "Sheet1"): new ID/Chinese column values + 5 extra
#    blank-but-formatted rows (8-12); wider columns.
# ---------------------------------------------------------------------------
$sheetComm.Range("A2").Value = "Langage_Comm_1"
$sheetComm.Range("C2").Value = "确认"

$sheetComm.Range("A3").Value = "Langage_Comm_2"
$sheetComm.Range("C3").Value = "取消"

$sheetComm.Range("A4").Value = "Langage_Comm_3"
$sheetComm.Range("C4").Value = "登录"

$sheetComm.Range("A5").Value = "Langage_Comm_4"
$sheetComm.Range("C5").Value = "创建角色"

$sheetComm.Range("A6").Value = "Langage_Comm_5"
$sheetComm.Range("C6").Value = "进入游戏"

$sheetComm.Range("A7").Value = "Langage_Comm_6"
$sheetComm.Range("C7").Value = "中文_6"

# Extend formatting (text number format + font) down to rows 8-12, matching
# the style already used by A2:C7.
$sheetComm.Range("A7:C7").Copy()
$sheetComm.Range("A8:C12").PasteSpecial(-4122)

$sheetComm.Columns.Item(1).ColumnWidth = 31.160714285714285
$sheetComm.Columns.Item(2).ColumnWidth = 23.785714285714285
$sheetComm.Columns.Item(3).ColumnWidth = 22.285714285714285

# ---------------------------------------------------------------------------
# 3. "Property" sheet: single formatted column of attribute-name labels.
# ---------------------------------------------------------------------------
$sheetProperty.Range("A1").Value = "ID"
$sheetProperty.Range("B1").Value = "English"
$sheetProperty.Range("C1").Value = "Chinese"

$sheetProperty.Range("A2").Value = "Langage_HP"
$sheetProperty.Range("A3").Value = "Langage_MAXHP"
$sheetProperty.Range("A4").Value = "Langage_MP"
$sheetProperty.Range("A5").Value = "Langage_MAXMP"
$sheetProperty.Range("A6").Value = "Langage_VP"
$sheetProperty.Range("A7").Value = "Langage_ATTACK"

# Give row-1 header cells (B1/C1) the same font-only style used in Comm.
$sheetComm.Range("B1:C1").Copy()
$sheetProperty.Range("B1:C1").PasteSpecial(-4122)

# A2's font-only style (no text numfmt here - matches authored file) extended
# down through the rest of the formatted-but-empty rows (8-28).
$sheetProperty.Range("A2").Copy()
$sheetProperty.Range("A3:A7").PasteSpecial(-4122)
$sheetProperty.Range("A8:A28").PasteSpecial(-4122)

$sheetProperty.Columns.Item(1).ColumnWidth = 50.535714285714285

# ---------------------------------------------------------------------------
# 4. "Guild" sheet: confirmation dialog strings + stray leftover formatting
#    (rows 3-12, 16 and 22) exactly as authored.
# ---------------------------------------------------------------------------
$sheetGuild.Range("A1").Value = "ID"
$sheetGuild.Range("B1").Value = "English"
$sheetGuild.Range("C1").Value = "Chinese"

$sheetGuild.Range("A2").Value = "Langage_Guild_1"
$sheetGuild.Range("B2").Value = "Langage_1"
$sheetGuild.Range("C2").Value = "确认要加入这个公会吗？点击确认加入"

$sheetGuild.Range("A2:C2").Copy()
$sheetGuild.Range("A3:C12").PasteSpecial(-4122)
$sheetGuild.Range("A16:C16").PasteSpecial(-4122)

$sheetGuild.Range("B1").Copy()
$sheetGuild.Range("A22").PasteSpecial(-4122)

$sheetGuild.Columns.Item(1).ColumnWidth = 31.160714285714285
$sheetGuild.Columns.Item(2).ColumnWidth = 23.785714285714285
$sheetGuild.Columns.Item(3).ColumnWidth = 22.285714285714285

# ---------------------------------------------------------------------------
# 5. "Tip" / "Item" sheets: header row only, nothing else yet.
# ---------------------------------------------------------------------------
$sheetTip.Range("A1").Value = "ID"
$sheetTip.Range("B1").Value = "English"
$sheetTip.Range("C1").Value = "Chinese"
$sheetComm.Range("B1:C1").Copy()
$sheetTip.Range("B1:C1").PasteSpecial(-4122)

$sheetItem.Range("A1").Value = "ID"
$sheetItem.Range("B1").Value = "English"
$sheetItem.Range("C1").Value = "Chinese"
$sheetComm.Range("B1:C1").Copy()
$sheetItem.Range("B1:C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Selections on each sheet, matching the authored file; Comm is done last
#    so it remains the active tab.
# ---------------------------------------------------------------------------
$sheetProperty.Range("A1:XFD1").Select()
$sheetGuild.Range("A12").Select()
$sheetTip.Range("A1:XFD1").Select()
$sheetItem.Range("A1:XFD1").Select()
$sheetComm.Range("C8").Select()
$sheetComm.Activate()
